$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45661 -> 45662) for every data row (rows 2-37).
foreach ($row in 2..37) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45661) {
        $cell.Value = 45662
    }
}
